# Fix input data reference errors
#
# The "BCpUC" sheet's B2 formula looked up year 2019 against the
# BBoSCpUC sheet's year header row, which only starts at 2020 - so the
# 2019 row always evaluated to #N/A. Delete that stale/erroring row;
# Excel shifts the remaining year rows (2020-2050) up by one and the
# formulas (which reference cells on their own row) re-resolve cleanly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCpUC")

$ws.Rows.Item(2).Delete()

# Leave the sheet in the same place the user left it after editing.
$ws.Activate()
$ws.Range("D12").Select()
